$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = 111814591
$ws.Range("B16").Value = 77515
$ws.Range("E16").Value = 6425
$ws.Range("F16").Value = "Garnlav"
$ws.Range("G16").Value = "Alectoria sarmentosa"
$ws.Range("H16").Value = "(Ach.) Ach."
$ws.Range("P16").Value = "åsele 1:1 (åsele 1:1), Ås lm"
$ws.Range("Q16").Value = 610012.4812897337
$ws.Range("R16").Value = 7121464.398116477
$ws.Range("S16").Value = 1
$ws.Range("Z16").Value = "17:50"
$ws.Range("AB16").Value = "17:50"
$ws.Range("A17").Value = 111815024
$ws.Range("B17").Value = 56414
$ws.Range("E17").Value = 100049
$ws.Range("F17").Value = "Spillkråka"
$ws.Range("G17").Value = "Dryocopus martius"
$ws.Range("H17").Value = "(Linnaeus, 1758)"
$ws.Range("P17").Value = "åsele 1:1 (åsele 1:1), Ås lm"
$ws.Range("Q17").Value = 609922.1399673244
$ws.Range("R17").Value = 7121488.212810148
$ws.Range("S17").Value = 1
$ws.Range("Z17").Value = "18:12"
$ws.Range("AB17").Value = "18:12"
$ws.Range("A18").Value = 111815269
$ws.Range("B18").Value = 90666
$ws.Range("D18").Value = "LC"
$ws.Range("E18").Value = 4364
$ws.Range("F18").Value = "Dropptaggsvamp"
$ws.Range("G18").Value = "Hydnellum ferrugineum"
$ws.Range("H18").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q18").Value = 610053.7842541422
$ws.Range("R18").Value = 7121273.15248157
$ws.Range("Z18").Value = "18:27"
$ws.Range("AB18").Value = "18:27"
$ws.Range("A19").Value = 111815114
$ws.Range("B19").Value = 90660
$ws.Range("E19").Value = 4362
$ws.Range("F19").Value = "Blå taggsvamp"
$ws.Range("G19").Value = "Hydnellum caeruleum"
$ws.Range("H19").Value = "(Hornem.) P.Karst."
$ws.Range("P19").Value = "åsele 1:1, Ås lm"
$ws.Range("Q19").Value = 610384.0265214761
$ws.Range("R19").Value = 7121170.261031131
$ws.Range("S19").Value = 5
$ws.Range("Z19").Value = "18:19"
$ws.Range("AB19").Value = "18:19"
$ws.Range("A20").Value = 111814478
$ws.Range("B20").Value = 77515
$ws.Range("D20").Value = "NT"
$ws.Range("E20").Value = 6425
$ws.Range("F20").Value = "Garnlav"
$ws.Range("G20").Value = "Alectoria sarmentosa"
$ws.Range("H20").Value = "(Ach.) Ach."
$ws.Range("Q20").Value = 610155.3487898401
$ws.Range("R20").Value = 7121461.207019502
$ws.Range("Z20").Value = "17:41"
$ws.Range("AB20").Value = "17:41"
$ws.Range("A21").Value = 111814688
$ws.Range("B21").Value = 90087
$ws.Range("D21").Value = "LC"
$ws.Range("E21").Value = 3298
$ws.Range("F21").Value = "Trådticka"
$ws.Range("G21").Value = "Climacocystis borealis"
$ws.Range("H21").Value = "(Fr.) Kotl. & Pouzar"
$ws.Range("Q21").Value = 610011.2059644217
$ws.Range("R21").Value = 7121475.688616944
$ws.Range("Z21").Value = "17:55"
$ws.Range("AB21").Value = "17:55"
$ws.Range("A22").Value = 111814104
$ws.Range("B22").Value = 56398
$ws.Range("E22").Value = 100109
$ws.Range("F22").Value = "Tretåig hackspett"
$ws.Range("G22").Value = "Picoides tridactylus"
$ws.Range("Q22").Value = 610154.5078508666
$ws.Range("R22").Value = 7121460.305022033
$ws.Range("Z22").Value = "17:23"
$ws.Range("AB22").Value = "17:23"
$ws.Range("A23").Value = 111814925
$ws.Range("B23").Value = 89686
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 658
$ws.Range("F23").Value = "Rosenticka"
$ws.Range("G23").Value = "Rhodofomes roseus"
$ws.Range("H23").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("P23").Value = "åsele 1:1, Ås lm"
$ws.Range("Q23").Value = 610384.0265214761
$ws.Range("R23").Value = 7121170.261031131
$ws.Range("S23").Value = 5
$ws.Range("Z23").Value = "18:08"
$ws.Range("AB23").Value = "18:08"
